$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2,3,4) are being cyclically rotated:
#   new row 2 <- old row 4
#   new row 3 <- old row 2
#   new row 4 <- old row 3
# Only columns A,B,E,F,G,H,Q,R,Z,AB actually change value (other columns
# happen to be identical across the three rows already).

$cols = @("A","B","E","F","G","H","Q","R","Z","AB")

# Capture current values for rows 2-4 before overwriting anything.
$rowData = @{}
foreach ($r in 2..4) {
    $rowData[$r] = @{}
    foreach ($col in $cols) {
        $rowData[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# Determine the new row order.
$newOrder = @{ 2 = 4; 3 = 2; 4 = 3 }

foreach ($r in 2..4) {
    $srcRow = $newOrder[$r]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $rowData[$srcRow][$col]
    }
}
